# Apply the "Integrate Agora data" edits to the BDPbES workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "About" ---
$about = $wb.Worksheets.Item("About")

# Update the explanatory notes text (rows 7-8), then remove the old row 9
# (its text is being dropped entirely as the methodology note was rewritten
# for the EU merit-order approach instead of the US priority-1 approach).
$about.Range("A7").Value = "We assign priority 2 to all of them as the merit-order is in place in the EU 28. "
$about.Range("A8").Value = "Even though there are support schemes for certain types we assume market-based dispatch according to least marginal cost."
$about.Rows.Item(9).Delete()

# --- Sheet "BDPbES" ---
$bdpbes = $wb.Worksheets.Item("BDPbES")

# petroleum (row 11) and natural gas peaker (row 12) no longer get a
# special priority-1 dispatch; they now share priority 2 with everything
# else, consistent with the new EU merit-order note above.
$bdpbes.Range("B11").Value = 2
$bdpbes.Range("B12").Value = 2

$wb.Save()
